$wb = $excel.ActiveWorkbook

# --- Sheet "ODI Batting": clear the stray empty cell B2 ---
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B2").Value = $null

# --- Add a new sheet "ODI Batting Extra" right after "ODI Bowling" ---
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsExtra = $wb.Worksheets.Add($null, $wsBowling)
$wsExtra.Name = "ODI Batting Extra"

# Header row (bold, centered, thin-bordered - matches the other sheets' header style)
$hdr = $wsExtra.Range("A1:F1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"

# Data rows (leading "'" forces text storage for numeric-looking values,
# matching the source data which stores these as text/inline strings)
$wsExtra.Cells.Item(2, 1).Value = "'4472"
$wsExtra.Cells.Item(2, 2).Value = 7
$wsExtra.Cells.Item(2, 3).Value = "'"
$wsExtra.Cells.Item(2, 4).Value = "'"
$wsExtra.Cells.Item(2, 5).Value = "'"
$wsExtra.Cells.Item(2, 6).Value = "NO"

$wsExtra.Cells.Item(3, 1).Value = "'4473"
$wsExtra.Cells.Item(3, 2).Value = 7
$wsExtra.Cells.Item(3, 3).Value = "'4"
$wsExtra.Cells.Item(3, 4).Value = "'0"
$wsExtra.Cells.Item(3, 5).Value = "'16.19%"
$wsExtra.Cells.Item(3, 6).Value = "YES"

$wsExtra.Cells.Item(4, 1).Value = "'4476"
$wsExtra.Cells.Item(4, 2).Value = 7
$wsExtra.Cells.Item(4, 3).Value = "'6"
$wsExtra.Cells.Item(4, 4).Value = "'3"
$wsExtra.Cells.Item(4, 5).Value = "'23.19%"
$wsExtra.Cells.Item(4, 6).Value = "NO"
